$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.228.67'
$ws.Range('E2').Value = '  +5.56%  '
$ws.Range('D3').Value = '2.594.56'
$ws.Range('E3').Value = '  +7.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '506.50'
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').Value = '156.23'
$ws.Range('E6').Value = '  +1.55%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -5.03%  '
$ws.Range('D9').Value = '2.628.45'
$ws.Range('E9').Value = '  +7.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.40'
$ws.Range('E10').Value = '  +3.19%  '
$ws.Range('E11').Value = '  +4.42%  '
$ws.Range('E12').Value = '  +3.00%  '
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = '3.051.46'
$ws.Range('E14').Value = '  +7.14%  '
$ws.Range('D15').Value = '60.535.91'
$ws.Range('E15').Value = '  +6.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.70'
$ws.Range('E16').Value = '  +5.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000140'
$ws.Range('E17').Value = '  +5.04%  '
$ws.Range('D18').Value = '2.621.72'
$ws.Range('E18').Value = '  +7.58%  '
$ws.Range('D19').Value = '4.78'
$ws.Range('E19').Value = '  +3.45%  '
$ws.Range('D20').Value = '343.62'
$ws.Range('E20').Value = '  +5.97%  '
$ws.Range('D21').Value = '10.41'
$ws.Range('E21').Value = '  +4.33%  '
$ws.Range('D22').Value = '6.16'
$ws.Range('E22').Value = '  +4.32%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').Value = '5.75'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '60.37'
$ws.Range('E25').Value = '  +4.60%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').Value = '0.423'
$ws.Range('E26').Value = '  +5.55%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.165'
$ws.Range('E27').Value = '  +3.45%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '0.992'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0852'
$ws.Range('E29').Value = '  +8.70%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '7.54'
$ws.Range('E30').Value = '  +4.01%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '156.31'
$ws.Range('E32').Value = '  +3.95%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '19.35'
$ws.Range('E33').Value = '  +3.65%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '1.57'
$ws.Range('E34').Value = '  +3.19%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.70'
$ws.Range('E35').Value = '  +7.88%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '4.01'
$ws.Range('E36').Value = '  +6.24%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.21'
$ws.Range('E37').Value = '  +5.83%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '310.38'
$ws.Range('E38').Value = '  +9.27%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.850'
$ws.Range('E39').Value = '  +3.74%  '
$ws.Range('D40').Value = '3.78'
$ws.Range('E40').Value = '  +7.31%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '1.47'
$ws.Range('E41').Value = '  +7.26%  '
$ws.Range('B42').Value = 'SuiNetwork'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D42').Value = '0.835'
$ws.Range('E42').Value = '  +27.72%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '35.69'
$ws.Range('E43').Value = '  +4.88%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0570'
$ws.Range('E44').Value = '  +7.32%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.625'
$ws.Range('E45').Value = '  +4.04%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.101'
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = '0.993'
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '19.89'
$ws.Range('E48').Value = '  +13.17%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '4.88'
$ws.Range('E49').Value = '  +6.90%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0236'
$ws.Range('E50').Value = '  +3.83%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.042.35'
$ws.Range('E51').Value = '  +7.80%  '
